$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 ("zp") values were shifted: old B13:D13 (ZP D columns) move to E13:G13 (ZP S columns)
# and old E13:G13 (ZP S columns) move to B13:D13 (ZP D columns), correcting the ZP D / ZP S
# inversion described in the commit message.
$ws.Range("B13").Value = 2705
$ws.Range("C13").Value = 2588
$ws.Range("D13").Value = 2708
$ws.Range("E13").Value = 2573
$ws.Range("F13").Value = 2577
$ws.Range("G13").Value = 2711
